# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 370
$wsExhibition.Range("F6").Value = 1979
$wsExhibition.Range("F7").Value = 105

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 370
$wsAll.Range("F10").Value = 1979
$wsAll.Range("F11").Value = 105
